$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value  = -8.109
$ws.Range("E5").Value  = 13.09

$ws.Range("D6").Value  = -7.788999999999999
$ws.Range("E6").Value  = 13.477

$ws.Range("C7").Value  = -12.529
$ws.Range("D7").Value  = -7.992

$ws.Range("A8").Value  = -21.007
$ws.Range("D8").Value  = -7.478

$ws.Range("D9").Value  = -7.828999999999999

$ws.Range("A10").Value = -20.864
$ws.Range("D10").Value = -7.478

$ws.Range("A12").Value = -21.88
$ws.Range("D12").Value = -8.364999999999998
$ws.Range("E12").Value = 13.049

$ws.Range("B13").Value = 6.502999999999998

$ws.Range("A18").Value = -21.649

$ws.Range("E19").Value = 12.988

$ws.Range("C20").Value = -13.391
$ws.Range("E20").Value = 13.164

$ws.Range("E23").Value = 13.171

$ws.Range("A25").Value = -21.938
$ws.Range("E25").Value = 12.992
